# Update RFCN Mini Drone Results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the "DataSet" column (old column D) without disturbing the
#    column-width metadata of the remaining columns: copy E:M left onto
#    D:L, then wipe the now-vacated column M.
# ------------------------------------------------------------------
$ws.Range("E1:M8").Copy()
$ws.Range("D1").PasteSpecial()
$ws.Columns("M").Clear()
$excel.CutCopyMode = $false

# Re-apply the (shifted) column widths explicitly so they match the
# post-edit layout.
$ws.Columns("A").ColumnWidth = 6.28515625
$ws.Columns("B").ColumnWidth = 10.5703125
$ws.Columns("C").ColumnWidth = 8.7109375
$ws.Columns("D").ColumnWidth = 15.42578125
$ws.Columns("E").ColumnWidth = 14.7109375
$ws.Columns("F").ColumnWidth = 15
$ws.Columns("G").ColumnWidth = 12.5703125
$ws.Columns("H").ColumnWidth = 13.28515625
$ws.Columns("I").ColumnWidth = 14.140625
$ws.Columns("J").ColumnWidth = 17
$ws.Columns("K").ColumnWidth = 9.85546875
$ws.Columns("L").ColumnWidth = 10.42578125

# ------------------------------------------------------------------
# 2. Header row (row 1) - values already shifted correctly by the
#    copy/paste above (they're plain shared-string header labels), so
#    nothing else is required there.
# ------------------------------------------------------------------

# ------------------------------------------------------------------
# 3. Fill in the new results row (row 2).
# ------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Atul Gupta"
$ws.Range("C2").Value = "RFCN"
$ws.Range("D2").Value = "Mini Drone"
$ws.Range("E2").Value = "Mini Drone"
$ws.Range("F2").Value = 100000
$ws.Range("G2").Value = 15347
$ws.Range("H2").Value = 26931
$ws.Range("I2").Value = 0
$ws.Range("J2").Formula = "=G2/(G2+H2)"
$ws.Range("K2").Formula = "=H2/(H2+I2)"
$ws.Range("L2").Formula = "=2*(J2*K2)/(J2+K2)"

# Borders + centered alignment on the raw count cells (F2:I2).
$ws.Range("F2:I2").Borders.LineStyle = 1
$ws.Range("F2:I2").HorizontalAlignment = -4108

# Percent format + borders + centered alignment on the computed ratios
# (J2:L2).
$ws.Range("J2:L2").Style = "Percent"
$ws.Range("J2:L2").Borders.LineStyle = 1
$ws.Range("J2:L2").HorizontalAlignment = -4108
$ws.Range("J2:L2").NumberFormat = "0%"

# ------------------------------------------------------------------
# 4. Selection / view state, matching the authored workbook.
# ------------------------------------------------------------------
$ws.Range("F2:L2").Select()
